# Replace the hard-coded "2022. Año del Quincentenario de Toluca, Capital
# del Estado de México". legend (split across three runs with spell-check
# proofErr markers around "Quincentenario") with a single templated run
# containing ${leyenda}.
#
# Build the search string via character codes to avoid any source-encoding
# ambiguity with the accented characters / curly quote used in the document.
$d = $word.ActiveDocument

$openQuote  = [char]0x201C   # "
$eAcute     = [char]0x00E9   # é
$nTilde     = [char]0x00F1   # ñ

$searchText = "$openQuote" + "2022. A" + "$nTilde" + "o del Quincentenario de Toluca, Capital del Estado de M" + "$eAcute" + "xico`"."
$replaceText = '${leyenda}'

$found = $d.Content.Find.Execute(
    $searchText,    # FindText
    $true,          # MatchCase
    $false,         # MatchWholeWord
    $false,         # MatchWildcards
    $false,         # MatchSoundsLike
    $false,         # MatchAllWordForms
    $true,          # Forward
    1,              # Wrap (wdFindContinue)
    $false,         # Format
    $replaceText,   # ReplaceWith
    2               # Replace (wdReplaceAll)
)

Write-Host "Leyenda placeholder replacement executed, found/replaced: $found"
